$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matches source formatting,
# avoids Excel auto-converting to numbers and dropping significant trailing zeros).

$ws.Range('D2').Value = '46.518.63'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '2.600.01'
$ws.Range('E3').Value = '  +6.99%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '308.05'
$ws.Range('E5').Value = '  +3.82%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '100.46'
$ws.Range('E6').Value = '  +2.58%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.604'
$ws.Range('E7').Value = '  +5.57%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.579'
$ws.Range('E9').Value = '  +12.95%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.38'
$ws.Range('E10').Value = '  +11.09%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0846'
$ws.Range('E11').Value = '  +7.16%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '54.45'
$ws.Range('E12').Value = '  +1.52%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '8.17'
$ws.Range('E13').Value = '  +13.87%  '
$ws.Range('D14').Value = '2.995.41'
$ws.Range('E14').Value = '  +6.90%  '
$ws.Range('E15').Value = '  +1.22%  '
$ws.Range('D16').Value = '2.601.20'
$ws.Range('E16').Value = '  +6.31%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.924'
$ws.Range('E17').Value = '  +8.45%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '15.03'
$ws.Range('E18').Value = '  +6.51%  '
$ws.Range('D19').Value = '46.597.95'
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('E20').Value = '  +6.80%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.06'
$ws.Range('E21').Value = '  +1.97%  '
$ws.Range('E22').Value = '  +8.43%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '277.53'
$ws.Range('E23').Value = '  +12.82%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '71.81'
$ws.Range('E24').Value = '  +6.16%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.05'
$ws.Range('E25').Value = '  +8.63%  '
$ws.Range('E26').Value = '  +10.67%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '29.27'
$ws.Range('E27').Value = '  +36.70%  '
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.67'
$ws.Range('E29').Value = '  +8.76%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.26'
$ws.Range('E30').Value = '  +1.48%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '39.10'
$ws.Range('E31').Value = '  -1.06%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.33'
$ws.Range('E32').Value = '  +13.33%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.60'
$ws.Range('E33').Value = '  -7.07%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.84'
$ws.Range('E34').Value = '  +3.16%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0841'
$ws.Range('E35').Value = '  +8.51%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.21'
$ws.Range('E36').Value = '  +10.00%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '151.77'
$ws.Range('E37').Value = '  +2.43%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.123'
$ws.Range('E38').Value = '  +8.46%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E39').Value = '  +5.75%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '23.13'
$ws.Range('E40').Value = '  +40.45%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '16.16'
$ws.Range('E41').Value = '  +4.96%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0333'
$ws.Range('E42').Value = '  +9.68%  '
$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.65'
$ws.Range('E43').Value = '  +11.48%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.09'
$ws.Range('E44').Value = '  +3.57%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.141.25'
$ws.Range('E45').Value = '  +8.22%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.998'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '93.14'
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.52'
$ws.Range('E48').Value = '  +9.97%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.79'
$ws.Range('E49').Value = '  -2.72%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '109.52'
$ws.Range('E50').Value = '  +8.01%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.201'
$ws.Range('E51').Value = '  +7.84%  '
